# Apply the authored changes to the workbook:
#   - "Measurement technique" -> "Technique"
#   - "Measured property" -> "Measured property #parameter"
#   - Update the remembered selection on each renamed sheet
#   - Leave "Measured property #parameter" as the active (selected) tab

$wb = $excel.ActiveWorkbook

$wsTechnique = $wb.Worksheets.Item("Measurement technique")
$wsTechnique.Name = "Technique"

$wsProperty = $wb.Worksheets.Item("Measured property")
$wsProperty.Name = "Measured property #parameter"

# Visit "Technique" and leave its selection at B114
$wsTechnique.Activate()
$wsTechnique.Range("B114").Select()

# Finish on "Measured property #parameter" with selection at C62 -
# this is the sheet left active/selected in the saved workbook
$wsProperty.Activate()
$wsProperty.Range("C62").Select()
